# Rebuild the LR-pair result table (Hgf -> Sdc2) with the refreshed TPM-based
# NATMI statistics: a Resolving-Mac target-cluster column is now included, so
# every Sending-cluster x Target-cluster combination (5 x 4) is recomputed and
# rewritten below, growing the sheet from 16 to 21 used rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numDataRows = 20
$numCols = 20
$data = New-Object 'object[,]' $numDataRows, $numCols

$row0 = @('ECs', 'Hgf', 'Sdc2', 'ECs', 1, 0.3333333333333333, 0.09364566666666667, 0.280937, 0.001628901849080777, 0.001628901849080777, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 0.1710499907806667, 1.539449917026, [double]"4.553388093797721E-05", [double]"4.553388093797721E-05")
for ($c = 0; $c -lt $numCols; $c++) { $data[0, $c] = $row0[$c] }
$row1 = @('ECs', 'Hgf', 'Sdc2', 'FAPs', 1, 0.3333333333333333, 0.09364566666666667, 0.280937, 0.001628901849080777, 0.001628901849080777, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 4.147572164191445, 37.328149477723, 0.001104092764016141, 0.001104092764016141)
for ($c = 0; $c -lt $numCols; $c++) { $data[1, $c] = $row1[$c] }
$row2 = @('ECs', 'Hgf', 'Sdc2', 'MuSCs', 1, 0.3333333333333333, 0.09364566666666667, 0.280937, 0.001628901849080777, 0.001628901849080777, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 1.789462807966222, 16.105165271696, 0.0004763589057736562, 0.0004763589057736562)
for ($c = 0; $c -lt $numCols; $c++) { $data[2, $c] = $row2[$c] }
$row3 = @('ECs', 'Hgf', 'Sdc2', 'Resolving-Mac', 1, 0.3333333333333333, 0.09364566666666667, 0.280937, 0.001628901849080777, 0.001628901849080777, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 0.01095520074544444, 0.09859680670899999, [double]"2.91629835300224E-06", [double]"2.91629835300224E-06")
for ($c = 0; $c -lt $numCols; $c++) { $data[3, $c] = $row3[$c] }
$row4 = @('FAPs', 'Hgf', 'Sdc2', 'ECs', 3, 1, 10.26742066666667, 30.802262, 0.1785947081647151, 0.178594708164715, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 18.75412149743066, 168.787093476876, 0.004992388081770575, 0.004992388081770575)
for ($c = 0; $c -lt $numCols; $c++) { $data[4, $c] = $row4[$c] }
$row5 = @('FAPs', 'Hgf', 'Sdc2', 'FAPs', 3, 1, 10.26742066666667, 30.802262, 0.1785947081647151, 0.178594708164715, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 454.7446739494331, 4092.702065544898, 0.1210540248864669, 0.1210540248864669)
for ($c = 0; $c -lt $numCols; $c++) { $data[5, $c] = $row5[$c] }
$row6 = @('FAPs', 'Hgf', 'Sdc2', 'MuSCs', 3, 1, 10.26742066666667, 30.802262, 0.1785947081647151, 0.178594708164715, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 196.1987999096996, 1765.789199187296, 0.05222854882651082, 0.05222854882651082)
for ($c = 0; $c -lt $numCols; $c++) { $data[6, $c] = $row6[$c] }
$row7 = @('FAPs', 'Hgf', 'Sdc2', 'Resolving-Mac', 3, 1, 10.26742066666667, 30.802262, 0.1785947081647151, 0.178594708164715, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 1.201141051637111, 10.810269464734, 0.000319746369966731, 0.0003197463699667309)
for ($c = 0; $c -lt $numCols; $c++) { $data[7, $c] = $row7[$c] }
$row8 = @('Inflammatory-Mac', 'Hgf', 'Sdc2', 'ECs', 3, 1, 23.67539566666666, 71.02618699999999, 0.4118171950916292, 0.4118171950916292, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 43.24467276128066, 389.202054851526, 0.01151182628965393, 0.01151182628965393)
for ($c = 0; $c -lt $numCols; $c++) { $data[8, $c] = $row8[$c] }
$row9 = @('Inflammatory-Mac', 'Hgf', 'Sdc2', 'FAPs', 3, 1, 23.67539566666666, 71.02618699999999, 0.4118171950916292, 0.4118171950916292, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 1048.584686708608, 9437.262180377473, 0.2791355326011075, 0.2791355326011074)
for ($c = 0; $c -lt $numCols; $c++) { $data[9, $c] = $row9[$c] }
$row10 = @('Inflammatory-Mac', 'Hgf', 'Sdc2', 'MuSCs', 3, 1, 23.67539566666666, 71.02618699999999, 0.4118171950916292, 0.4118171950916292, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 452.4100422092995, 4071.690379883696, 0.1204325408208783, 0.1204325408208783)
for ($c = 0; $c -lt $numCols; $c++) { $data[10, $c] = $row10[$c] }
$row11 = @('Inflammatory-Mac', 'Hgf', 'Sdc2', 'Resolving-Mac', 3, 1, 23.67539566666666, 71.02618699999999, 0.4118171950916292, 0.4118171950916292, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 2.76968194566211, 24.927137510959, 0.0007372953799895674, 0.0007372953799895674)
for ($c = 0; $c -lt $numCols; $c++) { $data[11, $c] = $row11[$c] }
$row12 = @('MuSCs', 'Hgf', 'Sdc2', 'ECs', 3, 1, 0.3314846666666666, 0.9944539999999999, 0.005765947381177186, 0.005765947381177185, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 0.6054786216546666, 5.449307594892, 0.0001611797308090255, 0.0001611797308090255)
for ($c = 0; $c -lt $numCols; $c++) { $data[12, $c] = $row12[$c] }
$row13 = @('MuSCs', 'Hgf', 'Sdc2', 'FAPs', 3, 1, 0.3314846666666666, 0.9944539999999999, 0.005765947381177186, 0.005765947381177185, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 14.68147566525178, 132.133280987266, 0.003908240870895994, 0.003908240870895993)
for ($c = 0; $c -lt $numCols; $c++) { $data[13, $c] = $row13[$c] }
$row14 = @('MuSCs', 'Hgf', 'Sdc2', 'MuSCs', 3, 1, 0.3314846666666666, 0.9944539999999999, 0.005765947381177186, 0.005765947381177185, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 6.334297181336889, 57.008674632032, 0.001686203737073563, 0.001686203737073563)
for ($c = 0; $c -lt $numCols; $c++) { $data[14, $c] = $row14[$c] }
$row15 = @('MuSCs', 'Hgf', 'Sdc2', 'Resolving-Mac', 3, 1, 0.3314846666666666, 0.9944539999999999, 0.005765947381177186, 0.005765947381177185, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 0.03877895471977777, 0.3490105924779999, [double]"1.032304239860356E-05", [double]"1.032304239860356E-05")
for ($c = 0; $c -lt $numCols; $c++) { $data[15, $c] = $row15[$c] }
$row16 = @('Resolving-Mac', 'Hgf', 'Sdc2', 'ECs', 3, 1, 23.12211433333333, 69.366343, 0.4021932475133977, 0.4021932475133977, 3, 1, 1.826566, 5.479698, 0.02795372904983374, 0.02795372904983374, 42.23406788937933, 380.106611004414, 0.01124280106666224, 0.01124280106666224)
for ($c = 0; $c -lt $numCols; $c++) { $data[16, $c] = $row16[$c] }
$row17 = @('Resolving-Mac', 'Hgf', 'Sdc2', 'FAPs', 3, 1, 23.12211433333333, 69.366343, 0.4021932475133977, 0.4021932475133977, 3, 1, 44.29005966666667, 132.870179, 0.6778141756295529, 0.6778141756295529, 1024.079823442822, 9216.718410985397, 0.2726122845070664, 0.2726122845070664)
for ($c = 0; $c -lt $numCols; $c++) { $data[17, $c] = $row17[$c] }
$row18 = @('Resolving-Mac', 'Hgf', 'Sdc2', 'MuSCs', 3, 1, 23.12211433333333, 69.366343, 0.4021932475133977, 0.4021932475133977, 3, 1, 19.10886933333333, 57.326608, 0.2924417490485847, 0.2924417490485847, 441.8374615060604, 3976.537153554544, 0.1176180967583484, 0.1176180967583484)
for ($c = 0; $c -lt $numCols; $c++) { $data[18, $c] = $row18[$c] }
$row19 = @('Resolving-Mac', 'Hgf', 'Sdc2', 'Resolving-Mac', 3, 1, 23.12211433333333, 69.366343, 0.4021932475133977, 0.4021932475133977, 2, 0.6666666666666666, 0.1169856666666667, 0.350957, 0.001790346272028586, 0.001790346272028586, 2.704955960027888, 24.344603640251, 0.0007200651813206821, 0.0007200651813206821)
for ($c = 0; $c -lt $numCols; $c++) { $data[19, $c] = $row19[$c] }

# Write the full A2:T21 block in one shot (mirrors Excel's paste-values behaviour)
$targetRange = $ws.Range("A2").Resize($numDataRows, $numCols)
$targetRange.Value2 = $data
